$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5018.875
$ws.Range("J40").Value = 5725
$ws.Range("L40").Value = 5725
$ws.Range("N40").Value = -6075

$ws.Range("H98").Value = 1454.0454
$ws.Range("I98").Value = 1199.45
$ws.Range("K98").Value = 1199.45
$ws.Range("M98").Value = 298.55

$ws.Range("H107").Value = 728.0909
$ws.Range("I107").Value = 728.0909
$ws.Range("K107").Value = 728.0909
$ws.Range("M107").Value = 1191.9091

$ws.Range("H112").Value = 3086.6875
$ws.Range("I112").Value = 1985
$ws.Range("J112").Value = 3160.1333
$ws.Range("K112").Value = 5955
$ws.Range("L112").Value = 9480.3999
$ws.Range("M112").Value = -4847
$ws.Range("N112").Value = -11696.3999

$ws.Range("H122").Value = 1454.0454
$ws.Range("I122").Value = 1199.45
$ws.Range("K122").Value = 3598.35
$ws.Range("M122").Value = -1148.35

$ws.Range("H129").Value = 1378
$ws.Range("I129").Value = 472.5
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 1417.5
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = 3582.5
$ws.Range("N129").Value = -25000

$ws.Range("H132").Value = 111132104
$ws.Range("I132").Value = 166680320
$ws.Range("K132").Value = 500040960
$ws.Range("M132").Value = -500038430

$ws.Range("H138").Value = 7834.477
$ws.Range("J138").Value = 8385.103999999999
$ws.Range("L138").Value = 25155.312
$ws.Range("N138").Value = -35435.312

$ws.Range("H141").Value = 3326.0833
$ws.Range("I141").Value = 2734.5557
$ws.Range("K141").Value = 8203.667099999999
$ws.Range("M141").Value = -3023.667099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15694.333
$ws.Range("I32").Value = 15694.333
$ws.Range("K32").Value = 15694.333
$ws.Range("M32").Value = -15407.333

$ws.Range("H45").Value = 1148
$ws.Range("I45").Value = 607.4286
$ws.Range("J45").Value = 3040
$ws.Range("K45").Value = 607.4286
$ws.Range("L45").Value = 3040
$ws.Range("M45").Value = -230.4286
$ws.Range("N45").Value = -3794

$ws.Range("H74").Value = 106903360
$ws.Range("I74").Value = 152716940
$ws.Range("K74").Value = 152716940
$ws.Range("M74").Value = -152716066

$ws.Range("H77").Value = 106903360
$ws.Range("I77").Value = 152716940
$ws.Range("K77").Value = 763584700
$ws.Range("M77").Value = -763580332

$ws.Range("H102").Value = 2916.1052
$ws.Range("I102").Value = 1372.2858
$ws.Range("K102").Value = 1372.2858
$ws.Range("M102").Value = 249.7141999999999

$ws.Range("H132").Value = 4840.3076
$ws.Range("I132").Value = 4819.522
$ws.Range("K132").Value = 14458.566
$ws.Range("M132").Value = -11928.566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2233.3333
$ws.Range("I20").Value = 1633.3334
$ws.Range("K20").Value = 1633.3334
$ws.Range("M20").Value = -1386.3334

$ws.Range("H106").Value = 22978.8
$ws.Range("J106").Value = 22978.8
$ws.Range("L106").Value = 22978.8
$ws.Range("N106").Value = -25502.8

$ws.Range("H134").Value = 5499.3335
$ws.Range("I134").Value = 3000
$ws.Range("J134").Value = 6749
$ws.Range("K134").Value = 9000
$ws.Range("L134").Value = 20247
$ws.Range("M134").Value = -6465
$ws.Range("N134").Value = -25317

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 2033
$ws.Range("I8").Value = 200
$ws.Range("J8").Value = 2949.5
$ws.Range("K8").Value = 200
$ws.Range("L8").Value = 2949.5
$ws.Range("M8").Value = -60
$ws.Range("N8").Value = -3229.5

$ws.Range("H28").Value = 17187.666
$ws.Range("J28").Value = 17187.666
$ws.Range("L28").Value = 17187.666
$ws.Range("N28").Value = -17677.666

$ws.Range("H31").Value = 3130.5386
$ws.Range("I31").Value = 3090.9092
$ws.Range("J31").Value = 3348.5
$ws.Range("K31").Value = 3090.9092
$ws.Range("L31").Value = 3348.5
$ws.Range("M31").Value = -2795.9092
$ws.Range("N31").Value = -3938.5

$ws.Range("H34").Value = 3130.5386
$ws.Range("I34").Value = 3090.9092
$ws.Range("J34").Value = 3348.5
$ws.Range("K34").Value = 3090.9092
$ws.Range("L34").Value = 3348.5
$ws.Range("M34").Value = -2888.9092
$ws.Range("N34").Value = -3752.5

$ws.Range("H57").Value = 971.4286
$ws.Range("I57").Value = 971.4286
$ws.Range("K57").Value = 971.4286
$ws.Range("M57").Value = -411.4286

$ws.Range("H58").Value = 4049.2144
$ws.Range("J58").Value = 4164.8
$ws.Range("L58").Value = 4164.8
$ws.Range("N58").Value = -4570.8

$ws.Range("H134").Value = 2674.2
$ws.Range("I134").Value = 2699.125
$ws.Range("J134").Value = 2574.5
$ws.Range("K134").Value = 8097.375
$ws.Range("L134").Value = 7723.5
$ws.Range("M134").Value = -5562.375
$ws.Range("N134").Value = -12793.5

$ws.Range("H136").Value = 4049.2144
$ws.Range("J136").Value = 4164.8
$ws.Range("L136").Value = 12494.4
$ws.Range("N136").Value = -17594.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 201197.4
$ws.Range("J34").Value = 333996.34
$ws.Range("L34").Value = 1001989.02
$ws.Range("N34").Value = -1002157.02

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H52").Value = 5038.1665
$ws.Range("J52").Value = 5038.1665
$ws.Range("L52").Value = 15114.4995
$ws.Range("N52").Value = -15646.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2999999
$ws.Range("J11").Value = 2999999
$ws.Range("L11").Value = 2999999
$ws.Range("N11").Value = -3000277

$ws.Range("H18").Value = 7001
$ws.Range("I18").Value = 7001
$ws.Range("K18").Value = 7001
$ws.Range("M18").Value = -6708

$ws.Range("H102").Value = 3154.2144
$ws.Range("I102").Value = 3140.4614
$ws.Range("K102").Value = 3140.4614
$ws.Range("M102").Value = -1518.4614

$ws.Range("H132").Value = 22224492
$ws.Range("I132").Value = 2418.5
$ws.Range("J132").Value = 111112780
$ws.Range("K132").Value = 7255.5
$ws.Range("L132").Value = 333338340
$ws.Range("M132").Value = -4725.5
$ws.Range("N132").Value = -333343400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 83341910
$ws.Range("I40").Value = 125007910
$ws.Range("J40").Value = 9898.5
$ws.Range("K40").Value = 125007910
$ws.Range("L40").Value = 9898.5
$ws.Range("M40").Value = -125007774
$ws.Range("N40").Value = -10170.5

$ws.Range("H53").Value = 9992
$ws.Range("I53").Value = 9992
$ws.Range("K53").Value = 9992
$ws.Range("M53").Value = -9474

$ws.Range("H68").Value = 3498.5
$ws.Range("I68").Value = 1350
$ws.Range("J68").Value = 5647
$ws.Range("K68").Value = 1350
$ws.Range("L68").Value = 5647
$ws.Range("M68").Value = -601
$ws.Range("N68").Value = -7145

$ws.Range("H71").Value = 3498.5
$ws.Range("I71").Value = 1350
$ws.Range("J71").Value = 5647
$ws.Range("K71").Value = 6750
$ws.Range("L71").Value = 28235
$ws.Range("M71").Value = -3006
$ws.Range("N71").Value = -35723

$ws.Range("H132").Value = 2902.7144
$ws.Range("I132").Value = 2902.7144
$ws.Range("K132").Value = 8708.143199999999
$ws.Range("M132").Value = -6178.143199999999

$ws.Range("H136").Value = 3175.8572
$ws.Range("I136").Value = 3130.3333
$ws.Range("J136").Value = 3449
$ws.Range("K136").Value = 9390.999899999999
$ws.Range("L136").Value = 10347
$ws.Range("M136").Value = -6840.999899999999
$ws.Range("N136").Value = -15447

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H70").Value = 50000
$ws.Range("J70").Value = 50000
$ws.Range("L70").Value = 50000
$ws.Range("N70").Value = -50630

$ws.Range("H73").Value = 50000
$ws.Range("J73").Value = 50000
$ws.Range("L73").Value = 50000
$ws.Range("N73").Value = -52184

$ws.Range("H108").Value = 20000
$ws.Range("J108").Value = 20000
$ws.Range("L108").Value = 20000
$ws.Range("N108").Value = -27680

$ws.Range("H122").Value = 6473.75
$ws.Range("I122").Value = 6473.75
$ws.Range("K122").Value = 19421.25
$ws.Range("M122").Value = -16971.25

$ws.Range("H123").Value = 63580
$ws.Range("J123").Value = 63580
$ws.Range("L123").Value = 63580
$ws.Range("N123").Value = -73380

$ws.Range("H132").Value = 333342660
$ws.Range("I132").Value = 20000
$ws.Range("J132").Value = 500004000
$ws.Range("K132").Value = 60000
$ws.Range("L132").Value = 1500012000
$ws.Range("M132").Value = -57470
$ws.Range("N132").Value = -1500017060

$ws.Range("H136").Value = 10040.471
$ws.Range("I136").Value = 12488
$ws.Range("K136").Value = 37464
$ws.Range("M136").Value = -34914
